# DA coordinates updated from Marja
# Fill in the DA low (DAL_14..DAL_20) coordinate rows that were
# previously blank, bump their row height slightly to fit, move the
# active selection to the newly-populated block, and renumber the
# A51:A53 "spacer" style now that the unused duplicate style was
# dropped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14 ---
$ws.Rows.Item(14).RowHeight = 13.4
$ws.Range("B14").Value = 40.435232
$ws.Range("C14").Value = -79.9137794596383

# --- Row 15 ---
$ws.Rows.Item(15).RowHeight = 13.4
$ws.Range("B15").Value = 40.434763
$ws.Range("C15").Value = -79.9137869977501

# --- Row 16 ---
$ws.Rows.Item(16).RowHeight = 13.4
$ws.Range("B16").Value = 40.434568
$ws.Range("C16").Value = -79.9140966771776

# --- Row 17 ---
$ws.Rows.Item(17).RowHeight = 13.4
$ws.Range("B17").Value = 40.434668
$ws.Range("C17").Value = -79.9151207932949

# --- Row 18 ---
$ws.Rows.Item(18).RowHeight = 13.4
$ws.Range("B18").Value = 40.434668
$ws.Range("C18").Value = -79.9151207932949

# --- Row 19 ---
$ws.Rows.Item(19).RowHeight = 13.4
$ws.Range("B19").Value = 40.435094
$ws.Range("C19").Value = -79.9163518926827

# --- Row 20 ---
$ws.Rows.Item(20).RowHeight = 13.4
$ws.Range("B20").Value = 40.435346
$ws.Range("C20").Value = -79.9163124785583

# The newly entered block lost the old "empty placeholder" styling
# (applyFont/applyAlignment) and reverts to the plain default style,
# i.e. not wrapped.
$ws.Range("B14:C20").WrapText = $false

# Move the active selection to the block we just filled in.
$ws.Range("B14:C20").Select()
